$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.527.66"
$ws.Range("E2").Value = "  -2.14%  "

# Row 3
$ws.Range("D3").Value = "3.786.66"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.48"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.52"
$ws.Range("E6").Value = "  -1.86%  "

# Row 7
$ws.Range("D7").Value = "3.787.02"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -1.38%  "

# Row 10
$ws.Range("E10").Value = "  -1.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.36"
$ws.Range("E11").Value = "  -2.19%  "

# Row 12
$ws.Range("E12").Value = "  -0.98%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  -2.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.98"
$ws.Range("E14").Value = "  -2.12%  "

# Row 15
$ws.Range("D15").Value = "4.417.74"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16
$ws.Range("D16").Value = "3.753.65"
$ws.Range("E16").Value = "  -1.09%  "

# Row 17
$ws.Range("D17").Value = "67.421.96"
$ws.Range("E17").Value = "  -2.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.14"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.99"
$ws.Range("E20").Value = "  -1.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  -6.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.81"
$ws.Range("E22").Value = "  -2.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.699"
$ws.Range("E23").Value = "  -1.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("E24").Value = "  +2.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.44"
$ws.Range("E25").Value = "  -1.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.14"
$ws.Range("E26").Value = "  -3.97%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.86"
$ws.Range("E27").Value = "  -2.67%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -1.88%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("E30").Value = "  -1.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.86"
$ws.Range("E31").Value = "  -1.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -1.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.19"
$ws.Range("E33").Value = "  -3.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.17"
$ws.Range("E34").Value = "  -2.30%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("D36").Value = "3.735.90"
$ws.Range("E36").Value = "  -0.45%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0998"
$ws.Range("E37").Value = "  -2.25%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  -5.53%  "

# Row 39
$ws.Range("E39").Value = "  -1.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.73"
$ws.Range("E41").Value = "  -2.33%  "

# Row 42
$ws.Range("E42").Value = "  -0.15%  "

# Row 43
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.07"
$ws.Range("E44").Value = "  +0.33%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").Value = "  -4.06%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.05"
$ws.Range("E46").Value = "  +1.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.36"
$ws.Range("E47").Value = "  -3.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.74"
$ws.Range("E48").Value = "  +0.99%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "392.43"
$ws.Range("E49").Value = "  -1.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.83"
$ws.Range("E50").Value = "  -7.57%  "

# Row 51
$ws.Range("D51").Value = "2.754.77"
$ws.Range("E51").Value = "  +1.98%  "
